# Apply updated cryptocurrency price/volume data to sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# need to be forced to Text format first, so they round-trip as strings
# (matching the source inline-string cell type) instead of being coerced
# into numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the new values
$ws.Range("D2").Value = '36.402.20'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '1.986.41'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '235.44'
$ws.Range("E5").Value = '  -9.25%  '
$ws.Range("E6").Value = '  -3.41%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '54.37'
$ws.Range("E8").Value = '  -3.11%  '
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("D10").Value = '58.01'
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("D11").Value = '0.0748'
$ws.Range("E11").Value = '  -3.49%  '
$ws.Range("D12").Value = '0.0986'
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("D13").Value = '14.20'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").Value = '2.281.04'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '20.15'
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("D16").Value = '0.754'
$ws.Range("E16").Value = '  -6.50%  '
$ws.Range("D17").Value = '5.04'
$ws.Range("E17").Value = '  -4.35%  '
$ws.Range("D18").Value = '1.988.52'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '36.377.18'
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").Value = '67.79'
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").Value = '0.0₃0802'
$ws.Range("E21").Value = '  -4.81%  '
$ws.Range("D22").Value = '5.26'
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("D23").Value = '221.07'
$ws.Range("E23").Value = '  -3.20%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  -9.95%  '
$ws.Range("D27").Value = '162.74'
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("D28").Value = '8.65'
$ws.Range("E28").Value = '  -3.91%  '
$ws.Range("D29").Value = '0.127'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").Value = '18.86'
$ws.Range("E30").Value = '  -4.30%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '0.116'
$ws.Range("E32").Value = '  -2.82%  '
$ws.Range("D33").Value = '4.37'
$ws.Range("E33").Value = '  -5.89%  '
$ws.Range("D34").Value = '0.0605'
$ws.Range("E34").Value = '  -6.85%  '
$ws.Range("D35").Value = '4.23'
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").Value = '2.33'
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -3.01%  '
$ws.Range("D39").Value = '3.31'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '5.47'
$ws.Range("E40").Value = '  +5.38%  '
$ws.Range("E41").Value = '  -1.48%  '
$ws.Range("D42").Value = '1.452.92'
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("E44").Value = '  -5.45%  '
$ws.Range("E45").Value = '  -9.67%  '
$ws.Range("D46").Value = '89.13'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '14.92'
$ws.Range("E47").Value = '  -4.93%  '
$ws.Range("D48").Value = '0.991'
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("E50").Value = '  -3.97%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '3.66'
$ws.Range("E51").Value = '  +6.73%  '
